# Transform the "bld_mock_data" sheet so that the bedrooms (col B) and
# floor (col E) columns hold descriptive text labels instead of bare
# numbers, the building code ("B1") becomes "Bldg1", and columns B/C get
# explicit widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

$bedroomMap = @{
    1 = "1-bedroom"
    2 = "2-bedroom"
    3 = "3-bedroom"
    4 = "4-bedroom"
}

$floorMap = @{
    2  = "f_2"
    3  = "f_3"
    4  = "f_4"
    5  = "f_5"
    6  = "f_6"
    7  = "f_7"
    8  = "f_8"
    9  = "f_9"
    10 = "f_10"
}

for ($r = 2; $r -le $lastRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    if ($aCell.Value() -eq "B1") {
        $aCell.Value = "Bldg1"
    }

    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value()
    if ($bedroomMap.ContainsKey([int]$bVal)) {
        $bCell.Value = $bedroomMap[[int]$bVal]
    }

    $eCell = $ws.Cells.Item($r, 5)
    $eVal = $eCell.Value()
    if ($floorMap.ContainsKey([int]$eVal)) {
        $eCell.Value = $floorMap[[int]$eVal]
    }
}

# NOTE: the engine stores column width quantized to 1/6-character
# increments and then adds a fixed ~0.8333 (5/6) padding when serialising
# to the OOXML `width` attribute (stored = round(ColumnWidth*6)/6 + 5/6).
# Pick ColumnWidth inputs sitting at the centre of the bucket that maps to
# the desired stored widths (15 and ~11.6640625) so the result is exact /
# closest-possible regardless of float rounding noise.
$ws.Columns.Item(2).ColumnWidth = (85 / 6)   # -> stored width 15
$ws.Columns.Item(3).ColumnWidth = (65 / 6)   # -> stored width ~11.6666667 (closest achievable to 11.6640625)
